$d = $word.ActiveDocument

function Replace-One($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
    }
    return $ok
}

# 1. Header row height for the first table row (cantSplit + tblHeader)
$t = $d.Tables.Item(1)
$r = $t.Rows.Item(1)
$r.Height = 36.4

# 2. "developed a number of methods" -> "developed methods"
Replace-One "developed a number of methods for training" "developed methods for training"

# 3. "labeld data" -> "labeled data"
Replace-One "My research demonstrated that labeld data can be obtained" "My research demonstrated that labeled data can be obtained"

# 4. "whie" -> "while"
Replace-One "require only a fraction of the data whie still exhibiting high performance" "require only a fraction of the data while still exhibiting high performance"

# 5. "large" -> "massive" savings
Replace-One "result in large savings" "result in massive savings"

# 6. "I lead the development" -> "I led the development"
Replace-One ", and other types of data mining. I lead the development of the methods" ", and other types of data mining. I led the development of the methods"

# 7. "as an" -> "as" (Loyola University Chicago as an assistant -> as assistant)
Replace-One "omputer science department at Loyola University Chicago as an assistant" "omputer science department at Loyola University Chicago as assistant"

# 8. insert "several " before "applications of deep neural networks" and "are" -> "were"
# (split into two smaller replacements so the literal straight quotes around "deep learning"
#  are never part of the replaced span, avoiding automatic smart-quote substitution)
Replace-One "focusing on applications of deep neural networks" "focusing on several applications of deep neural networks"
Replace-One "Our results are very positive: the neural models perform approximately" "Our results were very positive: the neural models perform approximately"

# 9. citation rework + "This work is directly relevant to" -> "These works are directly relevant to"
Replace-One "at least in several cases show considerable gains in performance, while completely eliminating the need for manual feature engineering. The manuscript describing our findings was just accepted for publication (Dligach, Miller, Lin, Bethard, and Savova, 2017). This work is directly relevant to this AHRQ proposal as a number of methods" "at least in several cases show considerable gains in performance, while completely eliminating the need for manual feature engineering (Dligach, Miller, Lin, Bethard, and Savova, 2017; Lin, Miller, Dligach, Bethard, and Savova, 2017). These works are directly relevant to this AHRQ proposal as a number of methods"

# 10. "a number of joint publications" -> "many joint publications"
Replace-One "been extremely productive and resulted in a number of joint publications." "been extremely productive and resulted in many joint publications."

# 11. insert new sentence about Dr. Afshar / Dr. Churpek collaboration
Replace-One "detecting alcohol abuse in trauma patients. Our team has on board leading experts in both NLP and quality measurement" "detecting alcohol abuse in trauma patients. Dr. Afshar, Dr. Churpek (University of Chicago site PI), and I are currently collaborating on cross-institutional validation of our NLP methods for ARDS detection. Our team has on board leading experts in both NLP and quality measurement"

# 12. remove "Accepted for publication in: " in references
Replace-One "Neural temporal information extraction. Accepted for publication in: Proceedings of the 15th Annual Meeting" "Neural temporal information extraction. Proceedings of the 15th Annual Meeting"

# 13. replace reference entry 3 (Miller TA, Bethard S, Dligach D, et al. Discovering narrative containers...) with new Lin et al. BioNLP 2017 reference
Replace-One "3 `tMiller TA, Bethard S, Dligach D, et al. Discovering narrative containers in clinical text. ACL 2013 2013;:18." "3 `tLin C, Miller TA, Dligach D, Bethard S, Savova GK. Representations of Time Expressions for Temporal Relation Extraction with Convolutional Neural Networks. In Proceedings of the 2017 Workshop on Biomedical Natural Language Processing (BioNLP 2017) held in conjunction with ACL-2017. August 2017. Vancouver, Canada."

# Re-bold "Dligach D" (author's own name) in the new reference, matching the convention used elsewhere.
$rngBold = $d.Content
$okBold = $rngBold.Find.Execute("Dligach D, Bethard S, Savova GK. Representations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($okBold) {
    $boldRng = $d.Range($rngBold.Start, $rngBold.Start + 9)
    $boldRng.Font.Bold = 1
} else {
    Write-Host "NOT FOUND for bolding: Dligach D, Bethard S, Savova GK. Representations"
}

# 14. trim the "This project focuses on..." paragraph (remove trailing content)
Replace-One "This project focuses on building a framework of open-source services that can be dynamically configured to transform EHR data into standards-conforming, comparable information suitable for large-scale analyses, inferencing, and integration of disparate health data. The clinical narrative and NLP methods for its processing are a central piece towards data normalization." "This project focuses on building a framework of open-source services that can be dynamically configured to transform EHR data into standards-conforming information suitable for large-scale analyses."

Write-Host "Done"
